$d = $word.ActiveDocument

# --- Step 1: "Prueba de documento JOM" -> "tercera" ---------------------
# A plain Find/Replace keeps the existing run (and the _GoBack bookmark
# that immediately follows it) intact.
$d.Content.Find.Execute(
    "Prueba de documento JOM", $true, $false, $false, $false, $false,
    $true, 1, $false, "tercera", 2) | Out-Null

# --- Step 2: find out exactly where "tercera" ends -----------------------
$locate = $d.Content
$locate.Find.Execute(
    "tercera", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$afterTercera = $locate.End

# --- Step 3: append " prueba de documento" as its own run, placed AFTER --
# --- the bookmarkStart/bookmarkEnd pair (i.e. after the original run)  --
#
# A plain InsertAfter on the collapsed end-of-paragraph range merges the
# new text into the previous run (same formatting) and also re-anchors
# the _GoBack bookmark past the freshly typed text. To reproduce the
# target markup -- two separate <w:r> elements with the bookmark sitting
# between them -- a one-character placeholder is typed right after
# "tercera" and its formatting is toggled, which forces a genuine run
# boundary (and leaves the bookmark anchored right after "tercera"). The
# real text is then inserted after that placeholder, which is finally
# deleted again.
$placeholder = $d.Range($afterTercera, $afterTercera)
$placeholder.InsertAfter("Z")
$placeholderMark = $d.Range($afterTercera, $afterTercera + 1)
$placeholderMark.Font.Bold = $true
$placeholderMark.Font.Bold = $false

$afterPlaceholder = $afterTercera + 1
$newTextRange = $d.Range($afterPlaceholder, $afterPlaceholder)
$newTextRange.InsertAfter(" prueba de documento")

$placeholderMark2 = $d.Range($afterTercera, $afterTercera + 1)
$placeholderMark2.Text = ""
